$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = 10000
$ws.Range("C14").Value = "petroleo"

$ws.Range("C14").Select()
